$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete column A (the "id" column) - remaining columns shift left,
# turning the old A1:H6 range into A1:G6.
$ws.Columns.Item(1).Delete() | Out-Null

# Mirror the author's final selection: the whole of column A selected
# (this is what triggered/preceded the column deletion in the UI).
$ws.Range("A1:A1048576").Select() | Out-Null
